$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 "Save" - copy formatting from existing header cell G1
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New column H values: 0 for rows 2-5, 1 for row 6
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
